$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 212.4076363333333
$ws.Range("H2").Value = 637.222909
$ws.Range("I2").Value = 0.5080632835800084
$ws.Range("J2").Value = 0.5080632835800084
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 13.441269
$ws.Range("N2").Value = 40.323807
$ws.Range("O2").Value = 0.08973082133481231
$ws.Range("P2").Value = 0.08973082133481232
$ws.Range("Q2").Value = 2855.028177610507
$ws.Range("R2").Value = 25695.25359849456
$ws.Range("S2").Value = 0.04558893572569581
$ws.Range("T2").Value = 0.04558893572569582
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 212.4076363333333
$ws.Range("H3").Value = 637.222909
$ws.Range("I3").Value = 0.5080632835800084
$ws.Range("J3").Value = 0.5080632835800084
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 54.711535
$ws.Range("N3").Value = 164.134605
$ws.Range("O3").Value = 0.3652416280068742
$ws.Range("P3").Value = 0.3652416280068742
$ws.Range("Q3").Value = 11621.14782951844
$ws.Range("R3").Value = 104590.3304656659
$ws.Range("S3").Value = 0.1855658608252805
$ws.Range("T3").Value = 0.1855658608252805
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 212.4076363333333
$ws.Range("H4").Value = 637.222909
$ws.Range("I4").Value = 0.5080632835800084
$ws.Range("J4").Value = 0.5080632835800084
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 63.67711
$ws.Range("N4").Value = 191.03133
$ws.Range("O4").Value = 0.4250937452800914
$ws.Range("P4").Value = 0.4250937452800915
$ws.Range("Q4").Value = 13525.50442363766
$ws.Range("R4").Value = 121729.539812739
$ws.Range("S4").Value = 0.2159745240563269
$ws.Range("T4").Value = 0.215974524056327
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 212.4076363333333
$ws.Range("H5").Value = 637.222909
$ws.Range("I5").Value = 0.5080632835800084
$ws.Range("J5").Value = 0.5080632835800084
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.96553866666667
$ws.Range("N5").Value = 53.896616
$ws.Range("O5").Value = 0.119933805378222
$ws.Range("P5").Value = 0.119933805378222
$ws.Range("Q5").Value = 3816.017603641772
$ws.Range("R5").Value = 34344.15843277594
$ws.Range("S5").Value = 0.06093396297270514
$ws.Range("T5").Value = 0.06093396297270515
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 161.9384456666667
$ws.Range("H6").Value = 485.815337
$ws.Range("I6").Value = 0.3873447295187379
$ws.Range("J6").Value = 0.3873447295187379
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.441269
$ws.Range("N6").Value = 40.323807
$ws.Range("O6").Value = 0.08973082133481231
$ws.Range("P6").Value = 0.08973082133481232
$ws.Range("Q6").Value = 2176.658209647551
$ws.Range("R6").Value = 19589.92388682796
$ws.Range("S6").Value = 0.03475676071942707
$ws.Range("T6").Value = 0.03475676071942708
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 161.9384456666667
$ws.Range("H7").Value = 485.815337
$ws.Range("I7").Value = 0.3873447295187379
$ws.Range("J7").Value = 0.3873447295187379
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 54.711535
$ws.Range("N7").Value = 164.134605
$ws.Range("O7").Value = 0.3652416280068742
$ws.Range("P7").Value = 0.3652416280068742
$ws.Range("Q7").Value = 8859.900937937431
$ws.Range("R7").Value = 79739.10844143688
$ws.Range("S7").Value = 0.1414744196093062
$ws.Range("T7").Value = 0.1414744196093062
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 161.9384456666667
$ws.Range("H8").Value = 485.815337
$ws.Range("I8").Value = 0.3873447295187379
$ws.Range("J8").Value = 0.3873447295187379
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 63.67711
$ws.Range("N8").Value = 191.03133
$ws.Range("O8").Value = 0.4250937452800914
$ws.Range("P8").Value = 0.4250937452800915
$ws.Range("Q8").Value = 10311.77221794536
$ws.Range("R8").Value = 92805.9499615082
$ws.Range("S8").Value = 0.1646578217856243
$ws.Range("T8").Value = 0.1646578217856243
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 161.9384456666667
$ws.Range("H9").Value = 485.815337
$ws.Range("I9").Value = 0.3873447295187379
$ws.Range("J9").Value = 0.3873447295187379
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.96553866666667
$ws.Range("N9").Value = 53.896616
$ws.Range("O9").Value = 0.119933805378222
$ws.Range("P9").Value = 0.119933805378222
$ws.Range("Q9").Value = 2909.311407244399
$ws.Range("R9").Value = 26183.80266519959
$ws.Range("S9").Value = 0.04645572740438036
$ws.Range("T9").Value = 0.04645572740438036
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.338549
$ws.Range("H10").Value = 1.015647
$ws.Range("I10").Value = 0.0008097840527861261
$ws.Range("J10").Value = 0.0008097840527861261
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.441269
$ws.Range("N10").Value = 40.323807
$ws.Range("O10").Value = 0.08973082133481231
$ws.Range("P10").Value = 0.08973082133481232
$ws.Range("Q10").Value = 4.550528178681
$ws.Range("R10").Value = 40.954753608129
$ws.Range("S10").Value = 0.0000726625881603321
$ws.Range("T10").Value = 0.00007266258816033211
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.338549
$ws.Range("H11").Value = 1.015647
$ws.Range("I11").Value = 0.0008097840527861261
$ws.Range("J11").Value = 0.0008097840527861261
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 54.711535
$ws.Range("N11").Value = 164.134605
$ws.Range("O11").Value = 0.3652416280068742
$ws.Range("P11").Value = 0.3652416280068742
$ws.Range("Q11").Value = 18.522535462715
$ws.Range("R11").Value = 166.702819164435
$ws.Range("S11").Value = 0.0002957668457736092
$ws.Range("T11").Value = 0.0002957668457736092
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.338549
$ws.Range("H12").Value = 1.015647
$ws.Range("I12").Value = 0.0008097840527861261
$ws.Range("J12").Value = 0.0008097840527861261
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 63.67711
$ws.Range("N12").Value = 191.03133
$ws.Range("O12").Value = 0.4250937452800914
$ws.Range("P12").Value = 0.4250937452800915
$ws.Range("Q12").Value = 21.55782191339
$ws.Range("R12").Value = 194.02039722051
$ws.Range("S12").Value = 0.0003442341358669456
$ws.Range("T12").Value = 0.0003442341358669456
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.338549
$ws.Range("H13").Value = 1.015647
$ws.Range("I13").Value = 0.0008097840527861261
$ws.Range("J13").Value = 0.0008097840527861261
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.96553866666667
$ws.Range("N13").Value = 53.896616
$ws.Range("O13").Value = 0.119933805378222
$ws.Range("P13").Value = 0.119933805378222
$ws.Range("Q13").Value = 6.082215150061334
$ws.Range("R13").Value = 54.739936350552
$ws.Range("S13").Value = 0.0000971204829852391
$ws.Range("T13").Value = 0.0000971204829852391
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 43.38855633333333
$ws.Range("H14").Value = 130.165669
$ws.Range("I14").Value = 0.1037822028484675
$ws.Range("J14").Value = 0.1037822028484675
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 13.441269
$ws.Range("N14").Value = 40.323807
$ws.Range("O14").Value = 0.08973082133481231
$ws.Range("P14").Value = 0.08973082133481232
$ws.Range("Q14").Value = 583.197257197987
$ws.Range("R14").Value = 5248.775314781884
$ws.Range("S14").Value = 0.009312462301529082
$ws.Range("T14").Value = 0.009312462301529083
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 43.38855633333333
$ws.Range("H15").Value = 130.165669
$ws.Range("I15").Value = 0.1037822028484675
$ws.Range("J15").Value = 0.1037822028484675
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 54.711535
$ws.Range("N15").Value = 164.134605
$ws.Range("O15").Value = 0.3652416280068742
$ws.Range("P15").Value = 0.3652416280068742
$ws.Range("Q15").Value = 2373.854518430638
$ws.Range("R15").Value = 21364.69066587575
$ws.Range("S15").Value = 0.03790558072651391
$ws.Range("T15").Value = 0.03790558072651391
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 43.38855633333333
$ws.Range("H16").Value = 130.165669
$ws.Range("I16").Value = 0.1037822028484675
$ws.Range("J16").Value = 0.1037822028484675
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 63.67711
$ws.Range("N16").Value = 191.03133
$ws.Range("O16").Value = 0.4250937452800914
$ws.Range("P16").Value = 0.4250937452800915
$ws.Range("Q16").Value = 2762.857874378863
$ws.Range("R16").Value = 24865.72086940977
$ws.Range("S16").Value = 0.0441171653022732
$ws.Range("T16").Value = 0.0441171653022732
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 43.38855633333333
$ws.Range("H17").Value = 130.165669
$ws.Range("I17").Value = 0.1037822028484675
$ws.Range("J17").Value = 0.1037822028484675
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.96553866666667
$ws.Range("N17").Value = 53.896616
$ws.Range("O17").Value = 0.119933805378222
$ws.Range("P17").Value = 0.119933805378222
$ws.Range("Q17").Value = 779.4987864973449
$ws.Range("R17").Value = 7015.489078476105
$ws.Range("S17").Value = 0.01244699451815125
$ws.Range("T17").Value = 0.01244699451815125
